$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: update D14
$ws.Range("D14").Value = 0.6447426901493167

# Row 15: update C15, add D15
$ws.Range("C15").Value = 0.2386249091493167
$ws.Range("D15").Value = 0.597740902

# Row 16: update B16, add C16
$ws.Range("B16").Value = -0.0107480648506833
$ws.Range("C16").Value = 0.042359665
